$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '27.607.60'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.647.90'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue 'D5' '212.66'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('E6').Value = '  +4.76%  '
$ws.Range('E7').Value = '  +0.01%  '
Set-TextValue 'D8' '23.56'
$ws.Range('E8').Value = '  -1.70%  '
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('E10').Value = '  -1.17%  '
Set-TextValue 'D11' '0.0888'
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('D12').Value = '1.881.03'
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').Value = '1.656.72'
$ws.Range('E13').Value = '  -0.47%  '
Set-TextValue 'D14' '0.585'
$ws.Range('E14').Value = '  +3.90%  '
$ws.Range('E15').Value = '  -2.46%  '
Set-TextValue 'D16' '64.48'
$ws.Range('E16').Value = '  -1.99%  '
$ws.Range('D17').Value = '27.561.28'
$ws.Range('E17').Value = '  +0.10%  '
Set-TextValue 'D18' '232.27'
$ws.Range('E18').Value = '  -3.44%  '
$ws.Range('E19').Value = '  -0.74%  '
Set-TextValue 'D20' '7.58'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('E22').Value = '  -3.09%  '
Set-TextValue 'D23' '9.75'
$ws.Range('E23').Value = '  +4.40%  '
$ws.Range('E24').Value = '  -1.50%  '
Set-TextValue 'D25' '149.12'
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('E26').Value = '  -2.58%  '
$ws.Range('E27').Value = '  +1.55%  '
$ws.Range('E28').Value = '  -0.03%  '
Set-TextValue 'D29' '15.60'
$ws.Range('E29').Value = '  -3.89%  '
$ws.Range('E30').Value = '  -1.99%  '
$ws.Range('E31').Value = '  -2.93%  '
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('E33').Value = '  +2.45%  '
$ws.Range('D34').Value = '1.425.73'
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('E35').Value = '  +2.84%  '
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  -3.97%  '
$ws.Range('E39').Value = '  -2.29%  '
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D42' '5.53'
$ws.Range('E42').Value = '  +2.49%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D43' '0.813'
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('E44').Value = '  +1.29%  '
Set-TextValue 'D45' '65.19'
$ws.Range('E45').Value = '  -5.62%  '
$ws.Range('D46').Value = '1.790.11'
$ws.Range('E47').Value = '  -1.38%  '
Set-TextValue 'D48' '88.27'
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('E49').Value = '  +1.14%  '
$ws.Range('E50').Value = '  -2.43%  '
$ws.Range('E51').Value = '  -0.10%  '
